# ---------------------------------------------------------------------------
# Report restructuring:
#   - "Tenencia" sheet gains extra instrument columns (RENTA VARIABLE / TASA
#     FIJA sub-headers, YMCQO/DGCU2/PBR/S31E5 tickers) with their amounts.
#   - The "Variacion - Retorno" sheet is renamed to "Referencias" and keeps
#     the reference-rate data (now labelled "Variables de Referencia").
#   - A brand new "Retorno" sheet is inserted in its place, mirroring the
#     "Tenencia" headers but holding return-rate figures.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- helper values -----------------------------------------------------
$valOnHardDollar = "ON HARD DOLLAR"
$valRentaVariable = "RENTA VARIABLE"
$valTasaFija = "TASA FIJA"
$valYmcqo = "YMCQO"
$valDgcu2 = "DGCU2"
$valPbr = "PBR"
$valS31E5 = "S31E5"
$valAmt1 = "222222222.64"
$valAmt2 = "9999999.64"
$valAmt3 = "1111111111.64"
$valAmt4 = "768873866.64"

# =========================================================================
# 1) "Tenencia" (sheet 1): add RENTA VARIABLE / TASA FIJA column groups
# =========================================================================
$wsTenencia = $wb.Worksheets.Item(1)

$wsTenencia.Range("B1").Value = $valOnHardDollar
$wsTenencia.Range("C1").Value = $valRentaVariable
$wsTenencia.Range("E1").Value = $valTasaFija
$wsTenencia.Range("C1:D1").Merge()
$wsTenencia.Range("E1:E1").Merge()

# Restore D1's original fill/font (merge copies C1's look onto it, but the
# target keeps D1 on its previous "B1-like" style) by pasting B1's format.
$wsTenencia.Range("B1").Copy()
$wsTenencia.Range("D1").PasteSpecial(-4122)
$wsTenencia.Application.CutCopyMode = $false

$wsTenencia.Range("B2").Value = $valYmcqo
$wsTenencia.Range("C2").Value = $valDgcu2
$wsTenencia.Range("D2").Value = $valPbr
$wsTenencia.Range("E2").Value = $valS31E5

# Match row-2's customFormat style (s="4") on the newly added cells.
$wsTenencia.Range("B2").Copy()
$wsTenencia.Range("C2:E2").PasteSpecial(-4122)
$wsTenencia.Application.CutCopyMode = $false

# Row 3/4/5 amounts - force plain text storage (these look numeric, and a
# bare Value= would silently convert them to floats / dates) then drop the
# temporary "@" number format so the cells stay styleless like the target.
$amountsRange = $wsTenencia.Range("B3:E5")
$amountsRange.NumberFormat = "@"

$wsTenencia.Range("B3").Value = $valAmt1
$wsTenencia.Range("C3").Value = $valAmt2
$wsTenencia.Range("D3").Value = $valAmt3
$wsTenencia.Range("E3").Value = $valAmt4

$wsTenencia.Range("B4").Value = $valAmt1
$wsTenencia.Range("C4").Value = $valAmt2
$wsTenencia.Range("D4").Value = $valAmt3
$wsTenencia.Range("E4").Value = $valAmt4

$wsTenencia.Range("B5").Value = $valAmt1
$wsTenencia.Range("C5").Value = $valAmt2
$wsTenencia.Range("D5").Value = $valAmt3
$wsTenencia.Range("E5").Value = $valAmt4

$amountsRange.ClearFormats()

# =========================================================================
# 2) Insert a new "Retorno" sheet right before the existing 2nd sheet, and
#    rename that existing 2nd sheet ("Variacion - Retorno") to
#    "Referencias". This keeps tab order Tenencia / Retorno / Referencias.
# =========================================================================
$wsOldSecond = $wb.Worksheets.Item(2)
$wsRetorno = $wb.Worksheets.Add($wsOldSecond)
$wsRetorno.Name = "Retorno"

# Same column widths / layout as Tenencia.
$wsRetorno.Columns.Item(1).ColumnWidth = $wsTenencia.Columns.Item(1).ColumnWidth

$wsRetorno.Range("B1").Value = $valOnHardDollar
$wsRetorno.Range("C1").Value = $valRentaVariable
$wsRetorno.Range("E1").Value = $valTasaFija
$wsRetorno.Range("B1:B1").Merge()
$wsRetorno.Range("C1:D1").Merge()
$wsRetorno.Range("E1:E1").Merge()

$wsTenencia.Range("B1").Copy()
$wsRetorno.Range("D1").PasteSpecial(-4122)
$wsTenencia.Range("C1").Copy()
$wsRetorno.Range("C1").PasteSpecial(-4122)
$wsTenencia.Range("E1").Copy()
$wsRetorno.Range("E1").PasteSpecial(-4122)
$wsRetorno.Application.CutCopyMode = $false

$wsRetorno.Range("A2").Value = "Fecha"
$wsRetorno.Range("B2").Value = $valYmcqo
$wsRetorno.Range("C2").Value = $valDgcu2
$wsRetorno.Range("D2").Value = $valPbr
$wsRetorno.Range("E2").Value = $valS31E5

$wsTenencia.Range("A2:E2").Copy()
$wsRetorno.Range("A2:E2").PasteSpecial(-4122)
$wsRetorno.Application.CutCopyMode = $false

$wsRetorno.Range("A3").Value = "2024-08-01"
$wsRetorno.Range("A4").Value = "2024-08-02"
$wsRetorno.Range("A5").Value = "2024-08-03"

$wsTenencia.Range("A3:A5").Copy()
$wsRetorno.Range("A3:A5").PasteSpecial(-4122)
$wsRetorno.Application.CutCopyMode = $false

$retornoBody = $wsRetorno.Range("B3:E5")
$retornoBody.NumberFormat = "@"

$wsRetorno.Range("B3").Value = "-"
$wsRetorno.Range("C3").Value = "-"
$wsRetorno.Range("D3").Value = "-"
$wsRetorno.Range("E3").Value = "-"

$wsRetorno.Range("B4").Value = "0.00"
$wsRetorno.Range("C4").Value = "0.00"
$wsRetorno.Range("D4").Value = "0.00"
$wsRetorno.Range("E4").Value = "0.00"

$wsRetorno.Range("B5").Value = "0.00"
$wsRetorno.Range("C5").Value = "0.00"
$wsRetorno.Range("D5").Value = "0.00"
$wsRetorno.Range("E5").Value = "0.00"

$retornoBody.ClearFormats()

# =========================================================================
# 3) Rename the original 2nd sheet to "Referencias" and refresh its labels.
# =========================================================================
$wsReferencias = $wb.Worksheets.Item(3)
$wsReferencias.Name = "Referencias"

$wsReferencias.Range("B1").Value = "Variables de Referencia"
$wsReferencias.Range("B2").Value = "Tasa de interés de préstamos por adelantos en cuenta corriente"

$refBody = $wsReferencias.Range("B3:B5")
$refBody.NumberFormat = "@"
$wsReferencias.Range("B3").Value = "0.00"
$wsReferencias.Range("B4").Value = "0.00"
$wsReferencias.Range("B5").Value = "0.00"
$refBody.ClearFormats()

$wsReferencias.Activate()

Write-Host "Report restructuring applied"
